$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Date the sheet was filled in (row 3, "Date")
$d = Get-Date -Year 2020 -Month 11 -Day 19 -Hour 0 -Minute 0 -Second 0
$ws.Range("B3").Value = $d.Date

# Team name and size
$ws.Range("B4").Value = "Limette"
$ws.Range("B5").Value = 4

# Team member names (column A) and their salaries (column B)
$ws.Range("A8").Value = "Lukas Hasler"
$ws.Range("B8").Value = 100

$ws.Range("A9").Value = "Pascal Strebel"
$ws.Range("B9").Value = 100

$ws.Range("A10").Value = "Cedric Weibel"
$ws.Range("B10").Value = 100

$ws.Range("A11").Value = "Robin Schmidiger"
$ws.Range("B11").Value = 100

# Row 12 ("Member 5") is no longer used - clear the placeholder label
$ws.Range("A12").Value = ""

# Tasks completed this week / to complete next week
# (shared-string append order follows a column-first fill: both "Tasks
# completed" entries, then both "Tasks to complete" entries)
$ws.Range("A19").Value = "Added more parts to the frontend"
$ws.Range("A20").Value = "Worked on our UI Library"

$ws.Range("B19").Value = "Finish development of the frontend."
$ws.Range("B20").Value = "Start building the backend."

# Leave the selection where the author's saved session left it
[void]$ws.Range("E6").Select()
